$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (cohort_year 2020, period_index 5): num_customers 29 -> 30
$ws.Range("C22").Value = 30
$ws.Range("E22").Value = 0.01130369253956292

# Row 27 (cohort_year 2021, period_index 4): num_customers 44 -> 45
$ws.Range("C27").Value = 45
$ws.Range("E27").Value = 0.01998223801065719

# Row 31 (cohort_year 2022, period_index 3): num_customers 45 -> 46
$ws.Range("C31").Value = 46
$ws.Range("E31").Value = 0.0198961937716263

# Row 34 (cohort_year 2023, period_index 2): num_customers 78 -> 80
$ws.Range("C34").Value = 80
$ws.Range("E34").Value = 0.03546099290780142

# Row 36 (cohort_year 2024, period_index 1): num_customers 130 -> 131
$ws.Range("C36").Value = 131
$ws.Range("E36").Value = 0.06787564766839378

# Row 37 (cohort_year 2025, period_index 0): num_customers 808 -> 817, cohort_size 808 -> 817
$ws.Range("C37").Value = 817
$ws.Range("D37").Value = 817
$ws.Range("E37").Value = 1
